$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 14:52"

# Row 4 (country index 8)
$ws.Range("B4").Value = 6175600
$ws.Range("C4").Value = 2364
$ws.Range("D4").Value = 3425907
$ws.Range("E4").Value = 2562461
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 187232

# Row 17 (country index 21)
$ws.Range("B17").Value = 315772
$ws.Range("C17").Value = 951
$ws.Range("D17").Value = 290796
$ws.Range("E17").Value = 21079
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 3897

# Row 45 (country index 49)
$ws.Range("B45").Value = 70667
$ws.Range("C45").Value = 527

# Row 81 (country index 85)
$ws.Range("B81").Value = 16985
$ws.Range("C81").Value = 94
$ws.Range("D81").Value = 15205
$ws.Range("E81").Value = 1156

# Row 85 (country index 89)
$ws.Range("B85").Value = 14341
$ws.Range("C85").Value = 11
$ws.Range("D85").Value = 11157
$ws.Range("E85").Value = 2581
$ws.Range("G85").Value = 3
$ws.Range("H85").Value = 603

# Row 141 (country index 145)
$ws.Range("B141").Value = 2107
$ws.Range("C141").Value = 2
$ws.Range("D141").Value = 1997
$ws.Range("E141").Value = 100
